$d = $word.ActiveDocument

# Start from the last paragraph currently in the document ("I am adding a doc file")
# and append the new paragraphs one at a time, mirroring what a user would do by
# placing the cursor at the end of that line and pressing Enter to type each new line.
$lines = @("Hey there", "Edit this", "hello")

foreach ($line in $lines) {
    $count = $d.Paragraphs.Count
    $r = $d.Paragraphs($count).Range
    $r.InsertParagraphAfter()
    $newR = $d.Paragraphs($count + 1).Range
    $newR.InsertBefore($line)
}
